$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new gamebook entry as row 35
$ws.Range("B35").Value = 1987
$ws.Range("C35").Value = "王子の対決"
$ws.Range("D35").Value = "Clash of Princes"
$ws.Range("E35").Value = "Shakaishisosha"
$ws.Range("F35").Value = "clash-of-princes.jpg"
$ws.Range("G35").Value = "slipcase set"

# Set custom width for column F (closest reachable value to the authored 25.6640625,
# this runtime quantizes ColumnWidth to 1/6-character increments)
$ws.Range("F1").ColumnWidth = 24.8333333333333

# Update the selection to match the post-edit state
$ws.Range("G36").Select()
